# Merge the three runs that make up "<id>p158r_1</id>" into a single run.
#
# Before:  [<id>] (Courier New, 7f6000, 9pt)  [p158r_1] (plain, black)  [</id>] (Courier New, 7f6000, 9pt)
# After:   [<id>p158r_1</id>] as one run, keeping the formatting of the
#          opening "<id>" run. Neighbouring runs (e.g. the following
#          "<head>" run) must stay untouched / unmerged.

$d = $word.ActiveDocument

# 1. Locate the "<id>" run -> gives us the formatting to keep and the
#    insertion point right after it.
$openRng = $d.Content
$foundOpen = $openRng.Find.Execute("<id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundOpen) {
    throw "could not find the '<id>' run"
}
$openEnd = $openRng.End

# 2. Locate the matching "</id>" run that follows it.
$afterOpen = $d.Range($openEnd, $d.Content.End)
$foundClose = $afterOpen.Find.Execute("</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundClose) {
    throw "could not find the '</id>' run"
}
$closeStart = $afterOpen.Start
$closeEnd = $afterOpen.End

# 3. Capture the text that needs to survive (the id value + closing tag)
#    before deleting anything.
$middleText = $d.Range($openEnd, $closeStart).Text
$closeTagText = $d.Range($closeStart, $closeEnd).Text

# 4. Delete right-to-left so earlier offsets stay valid: first the
#    "</id>" run, then the "p158r_1" run in between. Both runs disappear
#    entirely, leaving only the "<id>" run.
$closeRng = $d.Range($closeStart, $closeEnd)
$closeRng.Text = ""

$midRng = $d.Range($openEnd, $closeStart)
$midRng.Text = ""

# 5. Re-insert the combined text directly after the "<id>" run. Because
#    it abuts a run with identical formatting, Word folds it into that
#    same run (picking up Courier New / 7f6000 / 9pt), instead of
#    creating a new run or merging into whatever follows.
$insertionPoint = $d.Range($openEnd, $openEnd)
$insertionPoint.InsertAfter($middleText + $closeTagText)
